# The edit re-orders/reassigns the ink rows so that the Material (D) and
# Inventoryitem (H) values line up differently than before:
#   Row 2: PANTONE Yellow  -> Magenta
#   Row 3: Black           -> PANTONE Yellow (Inventoryitem cleared)
#   Row 5: Yellow          -> Black
#   Row 6: Magenta         -> Yellow
#   Row 7: Adhesive        -> Silicone
#   Row 8: Silicone        -> Adhesive
# (Row 4 - Cyan - and Row 9 - Roll - are unaffected.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Material")

$ws.Range("D2").Value = "Magenta - UV - "
$ws.Range("H2").Value = "10001836 - 9442 PRO MAGENTA BW5 UV - INK"

$ws.Range("D3").Value = "PANTONE Yellow U  - UV - "
$ws.Range("H3").Value = ""

$ws.Range("D5").Value = "Black - UV - "
$ws.Range("H5").Value = "10001817 - 9409 MIXING BLACK UV - INK"

$ws.Range("D6").Value = "Yellow - UV - "
$ws.Range("H6").Value = "10001305 - PROCESS YELLOW C UV"

$ws.Range("D7").Value = "Silicone"
$ws.Range("H7").Value = "10016451 - Evonik RW 10 Teco RC Silicone"

$ws.Range("D8").Value = "Adhesive"
$ws.Range("H8").Value = "10001053 - RAVENWOOD LINERLESS-ADHESIVE - 7445HD"
